# Fruta / hortaliza, semanal
# A new weekly price record is inserted at the top of the "Piña" data block
# (row 133), pushing the existing rows 133-231 down to 134-232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row before the current row 133. Excel shifts rows
# 133:231 down to 134:232 and the new blank row inherits formatting
# (including the date number-format in column D) from the row above.
$ws.Rows.Item(133).Insert()

# Populate the new row 133 with this week's record.
$ws.Cells.Item(133, 1).Value  = 5
$ws.Cells.Item(133, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value  = "Maule"
$ws.Cells.Item(133, 4).Value  = 44651
$ws.Cells.Item(133, 5).Value  = 7
$ws.Cells.Item(133, 6).Value  = "Fruta"
$ws.Cells.Item(133, 7).Value  = 100108
$ws.Cells.Item(133, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(133, 9).Value  = 100108005
$ws.Cells.Item(133, 10).Value = "Piña"
$ws.Cells.Item(133, 11).Value = "Caramelo"
$ws.Cells.Item(133, 12).Value = "Tercera"
$ws.Cells.Item(133, 13).Value = 230
$ws.Cells.Item(133, 14).Value = 16000
$ws.Cells.Item(133, 15).Value = 16000
$ws.Cells.Item(133, 16).Value = 16000
$ws.Cells.Item(133, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(133, 18).Value = "Ecuador"
$ws.Cells.Item(133, 19).Value = 1000
$ws.Cells.Item(133, 20).Value = 16
